# The document has three <id>...</id> tags that were each split across
# three separate runs (one run for the literal "<id>" text, one run for
# the bare id value, one run for the literal "</id>" text). Collapse each
# triple into a single run carrying the combined text, so the run picks
# up the formatting (Courier New / color 7f6000 / sz 18) of the first
# ("<id>") run, exactly as Word does when a Find/Replace match spans
# multiple runs.

$d = $word.ActiveDocument

$ids = @("p109r_5", "p109v_1", "p109v_2")

foreach ($id in $ids) {
    $tag = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($tag, $false, $false, $false, $false, $false, $true, 1, $false, $tag, 2)
}
